# The "Vovinam"/"Sao" certificate list rows were replaced by a proper
# subject-code/content table ("change list cert to table"): every data
# row now reuses the same SubjectCode ("ĐTR101" / the look-alike "ÐTR101"
# already used on row 2) while the Content column gets new song titles.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (student #2): SubjectCode stays ÐTR101, Content -> "Tình một đêm"
$ws.Range("I3").Value = "ÐTR101"
$ws.Range("K3").Value = "Tình một đêm"

# Row 4 (student #3): SubjectCode stays ÐTR101, Content -> "Sáo"
$ws.Range("I4").Value = "ÐTR101"
$ws.Range("K4").Value = "Sáo"

# Row 5 (student #4): SubjectCode -> "ĐTR101" (new code), Content -> "Lừa trái tim đàn bà"
$ws.Range("I5").Value = "ĐTR101"
$ws.Range("K5").Value = "Lừa trái tim đàn bà"

# Leave the selection on K3, as in the saved workbook.
$ws.Range("K3").Select()
